$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- C14: add SUM formula over column C ---
$ws.Range("C14").Formula = "=SUM(C3:C12)"

# --- Row 17: new task "set python selenium debugger" ---
$ws.Range("A17").Value = "set python selenium debugger"
$ws.Range("B17").Value = 1
$ws.Range("C17").Value = 0.75

# D17 needs the percentage style (same as D3:D5) -> copy format from D3
$ws.Range("D3").Copy() | Out-Null
$ws.Range("D17").PasteSpecial(-4122) | Out-Null
$ws.Range("D17").Value = 1

# --- Row 20: Final total row ---
$ws.Range("A20").HorizontalAlignment = -4108
$ws.Range("A20").VerticalAlignment = -4108
$ws.Range("A20").Value = "Final total"

$ws.Range("A20").Copy() | Out-Null
$ws.Range("B20:C20").PasteSpecial(-4122) | Out-Null

$ws.Range("B20").Formula = "=SUM(B14,B17:B18)"
$ws.Range("C20").Formula = "=SUM(C14, C17:C18)"

$excel.CutCopyMode = 0

# --- Update selection ---
$ws.Range("C15").Select() | Out-Null

Write-Output "done"
